# "models update by prakash"
#
# - Fix the "DegereePlanAbv" typo header (DegreePlan sheet, C1) to "DegreePlanAbv"
# - Fix the StudentID entry for DegreePlanID 7253 (DegreePlan!B4) to be the plain
#   numeric StudentID 534049 instead of the text "S534049"
# - Re-select cell C1 on the DegreePlan sheet (last-used cell after the edit)
# - Switch the active sheet from StudentTerm to Student

$wb = $excel.ActiveWorkbook

$degreePlan = $wb.Worksheets.Item("DegreePlan")
$degreePlan.Range("C1").Value = "DegreePlanAbv"
$degreePlan.Range("B4").Value = 534049
$degreePlan.Range("C1").Select()

$student = $wb.Worksheets.Item("Student")
$student.Activate()
